$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data region entirely before rewriting the final state,
# since row count / row content shifts extensively between versions.
$ws.Range("A1:C80").ClearContents()

$ws.Range("A1").Value = "Feature"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Description"

$ws.Range("B2").Value = "Interfaces"
$ws.Range("C2").Value = "Provide a GET interface"

$ws.Range("B3").Value = "Interfaces"
$ws.Range("C3").Value = "Provide a POST interface"

$ws.Range("C4").Value = "Method to Request streaming market data from IB Client"

$ws.Range("C5").Value = "Method to Request account information from IB Client"

$ws.Range("C6").Value = "External applications can subscribe to streaming market data via requests to the client"

$ws.Range("C7").Value = "External applications receive periodic updates to their market data subscriptions"

$ws.Range("C8").Value = "External applications shall be able to unsubscribe to market data"

$ws.Range("C9").Value = "External applications shall be able to submit, modify, and cancel orders"

$ws.Range("C10").Value = "External applications will receive updates about changes in order status (fills, etc.)"

$ws.Range("C11").Value = "External applications can request and receive status about the system"

$ws.Range("B12").Value = "Logging"
$ws.Range("C12").Value = "All loggable events will be stored in a master log file"

$ws.Range("B13").Value = "Logging"
$ws.Range("C13").Value = "Log files will cover a single day"

$ws.Range("B14").Value = "Logging"
$ws.Range("C14").Value = "Loggable events will be categorized in order to make event extraction easier"

$ws.Range("B15").Value = "Logging"
$ws.Range("C15").Value = "Loggable events will be timestamped"

$ws.Range("C16").Value = "Market data will be stored in a centralized database"

$ws.Range("B17").Value = "Interfaces"
$ws.Range("C17").Value = "provide a ROS2 interface"

$ws.Range("B19").Value = "Configuration"
$ws.Range("C19").Value = "Configuration parameters stored using ConfigObj format"

$ws.Range("B20").Value = "Configuration"
$ws.Range("C20").Value = "Configuration specification shall be supported"

$ws.Range("B21").Value = "Configuration"
$ws.Range("C21").Value = "configuration path shall default to root directory but optional path shall be supported"

$ws.Range("B22").Value = "Configuration"
$ws.Range("C22").Value = "parameter access shall be via a string path rather than multiple dictionaries"

$ws.Range("B23").Value = "Configuration"
$ws.Range("C23").Value = "A configured parameter may be optional or required"

$ws.Range("B25").Value = "User Interface"
$ws.Range("C25").Value = "The client will have a CLI"

$ws.Range("B26").Value = "Parameter"
$ws.Range("C26").Value = "Provide a configurable time between IB Client `"Tickle`" events"

$ws.Range("C27").Value = "CLI flags use standard Linux format"

$ws.Range("C28").Value = "System can recover and retrieve state from loss of network"

$ws.Range("C29").Value = "System can recover and retrieve state from host reboot"

$ws.Range("C30").Value = "System can recover and retrieve state from loss of communication with Client Portal"

$ws.Range("C31").Value = "Provide means to determine current version of Client Portal gateway"

$ws.Range("C32").Value = "Provide means to determine current version of library"

$ws.Range("C33").Value = "provide means to connect to ib websocket stream for market data"

# Resize the backing table (Table2) to cover the new data extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C76"))

# Restore selection to match the authored workbook state.
$ws.Range("C25").Select() | Out-Null
